# Horarios actualizados Linea 141 - 730
# Applies the 06:25:30 refresh: updated Hora_Scrap timestamps/Minutos for
# rows whose scrape time was superseded, a couple of re-ordered rows, and
# newly scraped rows appended at the bottom of each sheet.
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "Última actualización: 06:25:30"
$ws.Range("A3").Value = "Total filas: 56"
$ws.Range("A24").Value = "03:42:43"
$ws.Range("C24").Value = "14_ABASTO"
$ws.Range("D24").Value = 113
$ws.Range("A25").Value = "04:17:03"
$ws.Range("C25").Value = "215B_EL PATO"
$ws.Range("D25").Value = 78
$ws.Range("A36").Value = "06:25:30"
$ws.Range("D36").Value = 2
$ws.Range("A37").Value = "06:25:30"
$ws.Range("D37").Value = 4
$ws.Range("A38").Value = "06:25:30"
$ws.Range("D38").Value = 6
$ws.Range("A39").Value = "06:25:30"
$ws.Range("D39").Value = 19
$ws.Range("A40").Value = "06:25:30"
$ws.Range("D40").Value = 21
$ws.Range("A41").Value = "06:25:30"
$ws.Range("D41").Value = 34
$ws.Range("A42").Value = "06:25:30"
$ws.Range("B42").Value = "07:01"
$ws.Range("C42").Value = "16_SANTA ANA"
$ws.Range("D42").Value = 36
$ws.Range("A43").Value = "06:25:30"
$ws.Range("B43").Value = "07:04"
$ws.Range("C43").Value = "23_HERNANDEZ"
$ws.Range("D43").Value = 39
$ws.Range("A44").Value = "06:25:30"
$ws.Range("B44").Value = "07:05"
$ws.Range("C44").Value = "15_ABASTO"
$ws.Range("D44").Value = 40
$ws.Range("A45").Value = "05:27:50"
$ws.Range("B45").Value = "07:06"
$ws.Range("D45").Value = 99
$ws.Range("A46").Value = "06:25:30"
$ws.Range("B46").Value = "07:07"
$ws.Range("C46").Value = "225_GOMEZ"
$ws.Range("D46").Value = 42
$ws.Range("A47").Value = "06:25:30"
$ws.Range("B47").Value = "07:11"
$ws.Range("C47").Value = "215A_EL PATO"
$ws.Range("D47").Value = 46
$ws.Range("A48").Value = "06:25:30"
$ws.Range("B48").Value = "07:15"
$ws.Range("C48").Value = "11_ETCHEVERRY"
$ws.Range("D48").Value = 50
$ws.Range("A49").Value = "06:25:30"
$ws.Range("B49").Value = "07:21"
$ws.Range("C49").Value = "26_HERNANDEZ"
$ws.Range("D49").Value = 56
$ws.Range("A50").Value = "06:25:30"
$ws.Range("B50").Value = "07:23"
$ws.Range("C50").Value = "10_OLMOS"
$ws.Range("D50").Value = 58
$ws.Range("A52").Value = "06:25:30"
$ws.Range("B52").Value = "07:31"
$ws.Range("C52").Value = "11_ETCHEVERRY"
$ws.Range("D52").Value = 66
$ws.Range("A53").Value = "06:25:30"
$ws.Range("B53").Value = "07:32"
$ws.Range("C53").Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Range("D53").Value = 67
$ws.Range("A54").Value = "06:25:30"
$ws.Range("B54").Value = "07:36"
$ws.Range("C54").Value = "27_EL RETIRO"
$ws.Range("D54").Value = 71
$ws.Range("A55").Value = "06:25:30"
$ws.Range("B55").Value = "07:39"
$ws.Range("C55").Value = "10_OLMOS"
$ws.Range("D55").Value = 74
$ws.Range("A56").Value = "06:25:30"
$ws.Range("B56").Value = "07:47"
$ws.Range("C56").Value = "14_ABASTO"
$ws.Range("D56").Value = 82
$ws.Range("E56").Value = "LP1912"
$ws.Range("A57").Value = "06:25:30"
$ws.Range("B57").Value = "07:51"
$ws.Range("C57").Value = "215D_EL PATO"
$ws.Range("D57").Value = 86
$ws.Range("E57").Value = "LP1912"
$ws.Range("A58").Value = "06:25:30"
$ws.Range("B58").Value = "08:12"
$ws.Range("C58").Value = "15_ABASTO"
$ws.Range("D58").Value = 107
$ws.Range("E58").Value = "LP1912"
$ws.Range("A59").Value = "06:25:30"
$ws.Range("B59").Value = "08:21"
$ws.Range("C59").Value = "26_HERNANDEZ"
$ws.Range("D59").Value = 116
$ws.Range("E59").Value = "LP1912"
$ws.Range("A60").Value = "06:25:30"
$ws.Range("B60").Value = "08:22"
$ws.Range("C60").Value = "16_P MOR-SANTA ANA"
$ws.Range("D60").Value = 117
$ws.Range("E60").Value = "LP1912"
$ws.Range("A61").Value = "06:25:30"
$ws.Range("B61").Value = "08:23"
$ws.Range("C61").Value = "215B_EL PATO"
$ws.Range("D61").Value = 118
$ws.Range("E61").Value = "LP1912"

# ---- Sheet 2: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)

$ws.Range("A2").Value = "Última actualización: 06:25:30"
$ws.Range("A3").Value = "Total filas: 17"
$ws.Range("A19").Value = "06:25:30"
$ws.Range("D19").Value = 21
$ws.Range("A20").Value = "06:25:30"
$ws.Range("D20").Value = 46
$ws.Range("A21").Value = "06:25:30"
$ws.Range("D21").Value = 86
$ws.Range("A22").Value = "06:25:30"
$ws.Range("B22").Value = "08:23"
$ws.Range("C22").Value = "215B_EL PATO"
$ws.Range("D22").Value = 118
$ws.Range("E22").Value = "LP1912"

# ---- Sheet 3: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)

$ws.Range("A2").Value = "Última actualización: 06:25:30"
$ws.Range("A3").Value = "Total filas: 12"
$ws.Range("A12").Value = "06:25:30"
$ws.Range("D12").Value = 8
$ws.Range("A14").Value = "06:25:30"
$ws.Range("B14").Value = "07:00"
$ws.Range("C14").Value = "215B_LP-P MOR-1 Y 57"
$ws.Range("D14").Value = 35
$ws.Range("A15").Value = "05:55:46"
$ws.Range("B15").Value = "07:35"
$ws.Range("C15").Value = "215A_LA PLATA"
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = "L6173"
$ws.Range("A16").Value = "06:25:30"
$ws.Range("B16").Value = "07:40"
$ws.Range("C16").Value = "215A_LA PLATA"
$ws.Range("D16").Value = 75
$ws.Range("E16").Value = "L6173"
$ws.Range("A17").Value = "06:25:30"
$ws.Range("B17").Value = "08:07"
$ws.Range("C17").Value = "215C_LA PLATA"
$ws.Range("D17").Value = 102
$ws.Range("E17").Value = "L6203"

Write-Output "done"